# Edit script: remove stray gramStart/gramEnd proofErr markers around
# merged command runs, and restructure the "Push to Remote Server"
# section to add a note + hyperlink about adding an ssh key to GitHub.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Git config -global user.email ali@gmail.com" paragraph:
#    drop the gramStart/gramEnd proofErr wrapping "user.email" (keep
#    the spellStart/spellEnd pair and the mailto hyperlink untouched).
# ---------------------------------------------------------------------
$p97 = $d.Paragraphs.Item(97)
$xml97 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Git config –global </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>user.email</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ali@gmail.com</w:t></w:r></w:hyperlink></w:p>
'@
$p97.Range.InsertXML($xml97)
# InsertXML can't express rStyle directly, so reapply the Hyperlink
# character style to the e-mail address run afterwards.
$fix97 = $d.Content
[void]$fix97.Find.Execute("ali@gmail.com")
$fix97.Style = "Hyperlink"

# ---------------------------------------------------------------------
# 2) "Git add ." (first occurrence, "Adding in staging area" section):
#    merge the two split runs and drop gramStart/gramEnd.
# ---------------------------------------------------------------------
$p108 = $d.Paragraphs.Item(108)
$xml108 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Git add .</w:t></w:r><w:r><w:t xml:space="preserve"> -&gt; add all files inside current directory. But it has some draw backs when say directory structure is something like this:</w:t></w:r></w:p>
'@
$p108.Range.InsertXML($xml108)

# ---------------------------------------------------------------------
# 3) ".git" tree line: merge the empty tab run with the ".git" text run
#    and drop gramStart/gramEnd.
# ---------------------------------------------------------------------
$p110 = $d.Paragraphs.Item(110)
$xml110 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>.git</w:t></w:r></w:p>
'@
$p110.Range.InsertXML($xml110)

# ---------------------------------------------------------------------
# 4) "git add ." (second occurrence, test directory example): merge runs
#    and drop gramStart/gramEnd.
# ---------------------------------------------------------------------
$p115 = $d.Paragraphs.Item(115)
$xml115 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Now if you are in test directory and if you say </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>git add .</w:t></w:r><w:r><w:t xml:space="preserve"> , it will only add index.html and index.css and not index.js, so to add it in staging area you can use command: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>git add -A</w:t></w:r></w:p>
'@
$p115.Range.InsertXML($xml115)

# ---------------------------------------------------------------------
# 5) "git rm -cached -r ." merge runs and drop gramStart/gramEnd.
# ---------------------------------------------------------------------
$p118 = $d.Paragraphs.Item(118)
$xml118 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">To remove all files from staging area use command: </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>git rm –cached -r .</w:t></w:r></w:p>
'@
$p118.Range.InsertXML($xml118)

# ---------------------------------------------------------------------
# 6) "Push to Remote Server:" heading paragraph: split it into three
#    paragraphs -- the bare heading, a new note about needing an ssh
#    key, and the "adding origin:" sub-heading now preceded by a link
#    to GitHub's ssh docs.
# ---------------------------------------------------------------------
$p130 = $d.Paragraphs.Item(130)
$xml130 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:t>Push to Remote Server:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">*Note: You may need to add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ssh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-key in your </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> account, for pushing code to repo. You can follow steps here.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rStyle w:val="Heading5Char"/></w:rPr></w:pPr><w:r><w:t>SSHLINKPLACEHOLDER</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="Heading5Char"/></w:rPr><w:t>adding origin:</w:t></w:r></w:p>
'@
$p130.Range.InsertXML($xml130)

# Turn the placeholder text into a real hyperlink pointing at GitHub's
# ssh docs, with the URL itself as the visible text.
$linkRange = $d.Content
[void]$linkRange.Find.Execute("SSHLINKPLACEHOLDER")
[void]$d.Hyperlinks.Add($linkRange, "https://docs.github.com/en/authentication/connecting-to-github-with-ssh", $null, $null, "https://docs.github.com/en/authentication/connecting-to-github-with-ssh")

# InsertXML dropped the Heading5Char rStyle on the "adding origin:" run
# (same limitation as step 1); restore it explicitly.
$fixOrigin = $d.Content
[void]$fixOrigin.Find.Execute("adding origin:")
$fixOrigin.Style = "Heading5Char"

Write-Output "edit complete"
